# Chile coordinates + fix of csv delimiter
# Adds Region category (column E), Latitude (column G) and
# Longitude (column H) for the 16 Chilean subdivisions (rows 88-103)
# on sheet "iso3312". These values were missing because the source
# CSV used ";" as delimiter for these rows instead of ",".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iso3312")

$data = @(
    @{ Row = 88;  Lat = -46.035339;             Lon = -73.036434999999997 },
    @{ Row = 89;  Lat = -23.446308999999999;    Lon = -68.998489000000006 },
    @{ Row = 90;  Lat = -18.539467999999999;    Lon = -69.714607999999998 },
    @{ Row = 91;  Lat = -38.782228000000003;    Lon = -72.543251999999995 },
    @{ Row = 92;  Lat = -27.648325;             Lon = -70.432531999999995 },
    @{ Row = 93;  Lat = -37.225611999999998;    Lon = -73.108534000000006 },
    @{ Row = 94;  Lat = -30.540376999999999;    Lon = -70.967738999999995 },
    @{ Row = 95;  Lat = -34.434587999999998;    Lon = -71.154459000000003 },
    @{ Row = 96;  Lat = -42.046638000000002;    Lon = -73.008763999999999 },
    @{ Row = 97;  Lat = -39.922539;             Lon = -72.588984999999994 },
    @{ Row = 98;  Lat = -51.902416000000002;    Lon = -73.244017999999997 },
    @{ Row = 99;  Lat = -35.500971999999997;    Lon = -71.727129000000005 },
    @{ Row = 100; Lat = -36.602809999999998;    Lon = -72.073119000000005 },
    @{ Row = 101; Lat = -33.478729000000001;    Lon = -70.590025999999995 },
    @{ Row = 102; Lat = -20.099081000000002;    Lon = -69.456920999999994 },
    @{ Row = 103; Lat = -32.740869000000004;    Lon = -71.404539 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("E$r").Value() = "Region"
    $ws.Range("G$r").Value() = $item.Lat
    $ws.Range("H$r").Value() = $item.Lon
}

# Restore the selection left by the author at the end of editing.
$ws.Range("F101").Select()
